$wb = $excel.ActiveWorkbook

# --- Sheet 1: Login ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("I2").Value = "Fail"

# --- Sheet 2: Members ---
$ws2 = $wb.Worksheets.Item(2)

# Update existing row 2 (TC-0001)
$ws2.Range("G2").Value = "Active"
$ws2.Range("H2").Value = "Saved"
$ws2.Range("I2").Value = "Saved"

# Update F2 date format to custom yyyy/mm/dd
$ws2.Range("F2").NumberFormat = "yyyy/mm/dd"
$ws2.Columns.Item(6).NumberFormat = "yyyy/mm/dd"

# Add new row 3 (TC-0002)
$ws2.Range("A3").Value = "TC-0002"
$ws2.Range("B3").Value = "Test the Existing Member Email ID Should not be Allowed."
$ws2.Range("C3").Value = "Kathir"
$ws2.Range("D3").Value = "kathir.s@gmail.com"
$ws2.Range("E3").Value = 8956234578
$ws2.Range("F3").Value = 45962
$ws2.Range("F3").NumberFormat = "yyyy/mm/dd"
$ws2.Range("G3").Value = "Active"
$ws2.Range("H3").Value = "This is member already added."

# Hyperlink for D3
$ws2.Hyperlinks.Add($ws2.Range("D3"), "mailto:kathir.s@gmail.com")
$ws2.Range("D3").Style = "Hyperlink"

# Update selection to K10
$ws2.Range("K10").Select()

Write-Host "done"
